$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

function Clear-EmptyCell($row, $col) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.ClearFormats()
}

# Row 192
Set-TextCell 192 1 "2024-05-23"
Set-TextCell 192 2 "09:26:36"
Set-TextCell 192 3 "Palet atascado en la curva"
Set-TextCell 192 4 "-"
Set-TextCell 192 5 "-"
Set-TextCell 192 6 "-"
Set-TextCell 192 7 "-"
Set-TextCell 192 8 "09:26:40"
Set-TextCell 192 9 "0:00:04"

# Row 193
Set-TextCell 193 1 "2024-05-23"
Set-TextCell 193 2 "09:34:23"
Set-TextCell 193 3 "-"
Set-TextCell 193 4 "-"
Set-TextCell 193 5 "No coloca bien el sealling"
Set-TextCell 193 6 "-"
Set-TextCell 193 7 "-"
Set-TextCell 193 8 "09:34:25"
Set-TextCell 193 9 "0:00:02"

# Row 194
Set-TextCell 194 1 "2024-05-23"
Set-TextCell 194 2 "09:42:46"
Set-TextCell 194 3 "-"
Set-TextCell 194 4 "-"
Set-TextCell 194 5 "-"
Set-TextCell 194 6 "Robot no coloca bien filter en palet"
Set-TextCell 194 7 "-"
Set-TextCell 194 8 "09:42:48"
Set-TextCell 194 9 "0:00:02"

# Row 195
Set-TextCell 195 1 "2024-05-23"
Set-TextCell 195 2 "09:42:50"
Set-TextCell 195 3 "-"
Set-TextCell 195 4 "-"
Set-TextCell 195 5 "-"
Set-TextCell 195 6 "NOK Soldadura Plástico"
Set-TextCell 195 7 "-"
Set-TextCell 195 8 "09:42:52"
Set-TextCell 195 9 "0:00:02"

# Row 196 (repair time not yet recorded -> H/I left blank)
Set-TextCell 196 1 "2024-05-23"
Set-TextCell 196 2 "09:42:53"
Set-TextCell 196 3 "-"
Set-TextCell 196 4 "-"
Set-TextCell 196 5 "-"
Set-TextCell 196 6 "Pieza enganchada en HV Test"
Set-TextCell 196 7 "-"
Clear-EmptyCell 196 8
Clear-EmptyCell 196 9
